# Auto-generated edit script applying numeric cell updates per the commit diff.
# Source sheet: "Sheets/Ultros_Profits.xlsx" maps onto this workbook's 8 worksheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR); diff hunks were matched to sheets/rows
# by the unique combination of pre-existing cell values at each row.

$wb = $excel.ActiveWorkbook


# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H82").Value = 1058.5555
$ws.Range("I82").Value = 885.5294
$ws.Range("K82").Value = 2656.5882
$ws.Range("M82").Value = -2250.5882
$ws.Range("H85").Value = 1058.5555
$ws.Range("I85").Value = 885.5294
$ws.Range("K85").Value = 2656.5882
$ws.Range("M85").Value = -1252.5882
$ws.Range("H95").Value = 59000
$ws.Range("J95").Value = 59000
$ws.Range("L95").Value = 59000
$ws.Range("N95").Value = -64492
$ws.Range("H98").Value = 2389.5557
$ws.Range("I98").Value = 2478.96
$ws.Range("K98").Value = 2478.96
$ws.Range("M98").Value = -980.96
$ws.Range("H122").Value = 2389.5557
$ws.Range("I122").Value = 2478.96
$ws.Range("K122").Value = 7436.88
$ws.Range("M122").Value = -4986.88
$ws.Range("H132").Value = 43861.152
$ws.Range("I132").Value = 3731.6667
$ws.Range("K132").Value = 11195.0001
$ws.Range("M132").Value = -8665.000100000001
$ws.Range("H137").Value = 8011831
$ws.Range("I137").Value = 28607600
$ws.Range("K137").Value = 85822800
$ws.Range("M137").Value = -85820250
$ws.Range("H141").Value = 6873.355
$ws.Range("I141").Value = 4618.3076
$ws.Range("J141").Value = 18599.6
$ws.Range("K141").Value = 13854.9228
$ws.Range("L141").Value = 55798.8
$ws.Range("M141").Value = -8674.9228
$ws.Range("N141").Value = -66158.79999999999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 5052
$ws.Range("I45").Value = 3316.6667
$ws.Range("K45").Value = 3316.6667
$ws.Range("M45").Value = -2939.6667
$ws.Range("H61").Value = 4282.8667
$ws.Range("I61").Value = 3107.5652
$ws.Range("K61").Value = 3107.5652
$ws.Range("M61").Value = -2895.5652
$ws.Range("H63").Value = 12332.167
$ws.Range("I63").Value = 7999.3335
$ws.Range("J63").Value = 16665
$ws.Range("K63").Value = 7999.3335
$ws.Range("L63").Value = 16665
$ws.Range("M63").Value = -7313.3335
$ws.Range("N63").Value = -18037
$ws.Range("H66").Value = 12332.167
$ws.Range("I66").Value = 7999.3335
$ws.Range("J66").Value = 16665
$ws.Range("K66").Value = 39996.6675
$ws.Range("L66").Value = 83325
$ws.Range("M66").Value = -36564.6675
$ws.Range("N66").Value = -90189
$ws.Range("H136").Value = 4282.8667
$ws.Range("I136").Value = 3107.5652
$ws.Range("K136").Value = 9322.695599999999
$ws.Range("M136").Value = -6772.695599999999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 47390.5
$ws.Range("I86").Value = 68223
$ws.Range("K86").Value = 68223
$ws.Range("M86").Value = -67100
$ws.Range("H87").Value = 30000
$ws.Range("I87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("M87").ClearContents()
$ws.Range("H89").Value = 47390.5
$ws.Range("I89").Value = 68223
$ws.Range("K89").Value = 341115
$ws.Range("M89").Value = -335499
$ws.Range("H90").Value = 30000
$ws.Range("I90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("M90").ClearContents()
$ws.Range("H99").Value = 30755.875
$ws.Range("I99").Value = 42369.4
$ws.Range("K99").Value = 42369.4
$ws.Range("M99").Value = -40871.4
$ws.Range("H105").Value = 2164.8
$ws.Range("I105").Value = 2330.2258
$ws.Range("J105").Value = 882.75
$ws.Range("K105").Value = 2330.2258
$ws.Range("L105").Value = 882.75
$ws.Range("M105").Value = -583.2258000000002
$ws.Range("N105").Value = -4376.75

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6008.724
$ws.Range("I31").Value = 8226.8125
$ws.Range("J31").Value = 3278.7693
$ws.Range("K31").Value = 8226.8125
$ws.Range("L31").Value = 3278.7693
$ws.Range("M31").Value = -7931.8125
$ws.Range("N31").Value = -3868.7693
$ws.Range("H34").Value = 6008.724
$ws.Range("I34").Value = 8226.8125
$ws.Range("J34").Value = 3278.7693
$ws.Range("K34").Value = 8226.8125
$ws.Range("L34").Value = 3278.7693
$ws.Range("M34").Value = -8024.8125
$ws.Range("N34").Value = -3682.7693
$ws.Range("H50").Value = 0
$ws.Range("I50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("M50").ClearContents()
$ws.Range("H54").Value = 22000
$ws.Range("J54").Value = 22000
$ws.Range("L54").Value = 22000
$ws.Range("N54").Value = -23316
$ws.Range("H56").Value = 7546.5
$ws.Range("I56").Value = 93
$ws.Range("J56").Value = 15000
$ws.Range("K56").Value = 93
$ws.Range("L56").Value = 15000
$ws.Range("M56").Value = 752
$ws.Range("N56").Value = -16690
$ws.Range("H59").Value = 368564
$ws.Range("I59").Value = 368564
$ws.Range("K59").Value = 368564
$ws.Range("M59").Value = -367419
$ws.Range("H63").Value = 50000
$ws.Range("J63").Value = 50000
$ws.Range("L63").Value = 50000
$ws.Range("N63").Value = -51372
$ws.Range("H66").Value = 50000
$ws.Range("J66").Value = 50000
$ws.Range("L66").Value = 150000
$ws.Range("N66").Value = -156864
$ws.Range("H68").Value = 55000
$ws.Range("J68").Value = 55000
$ws.Range("L68").Value = 55000
$ws.Range("N68").Value = -56498
$ws.Range("H69").Value = 9133.166999999999
$ws.Range("I69").Value = 6959.8
$ws.Range("K69").Value = 6959.8
$ws.Range("M69").Value = -6210.8
$ws.Range("H71").Value = 55000
$ws.Range("J71").Value = 55000
$ws.Range("L71").Value = 165000
$ws.Range("N71").Value = -172488
$ws.Range("H72").Value = 9133.166999999999
$ws.Range("I72").Value = 6959.8
$ws.Range("K72").Value = 20879.4
$ws.Range("M72").Value = -17135.4
$ws.Range("H88").Value = 35000
$ws.Range("J88").Value = 35000
$ws.Range("L88").Value = 35000
$ws.Range("N88").Value = -35812
$ws.Range("H91").Value = 35000
$ws.Range("J91").Value = 35000
$ws.Range("L91").Value = 35000
$ws.Range("N91").Value = -37808
$ws.Range("H122").Value = 3590.4736
$ws.Range("I122").Value = 4790.3
$ws.Range("J122").Value = 2257.3333
$ws.Range("K122").Value = 14370.9
$ws.Range("L122").Value = 6771.999899999999
$ws.Range("M122").Value = -11920.9
$ws.Range("N122").Value = -11671.9999
$ws.Range("H132").Value = 2899.647
$ws.Range("I132").Value = 2353.1333
$ws.Range("K132").Value = 7059.3999
$ws.Range("M132").Value = -4529.3999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 193332.33
$ws.Range("J37").Value = 193332.33
$ws.Range("L37").Value = 579996.99
$ws.Range("N37").Value = -580220.99
$ws.Range("H109").Value = 499
$ws.Range("I109").Value = 499
$ws.Range("K109").Value = 1497
$ws.Range("M109").Value = -457
$ws.Range("H133").Value = 0
$ws.Range("I133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("M133").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4849.923
$ws.Range("I122").Value = 4368.136
$ws.Range("K122").Value = 13104.408
$ws.Range("M122").Value = -10654.408
$ws.Range("H124").Value = 39833.332
$ws.Range("J124").Value = 39833.332
$ws.Range("L124").Value = 39833.332
$ws.Range("N124").Value = -49653.332
$ws.Range("H126").Value = 3277
$ws.Range("I126").Value = 3043.75
$ws.Range("J126").Value = 3743.5
$ws.Range("K126").Value = 9131.25
$ws.Range("L126").Value = 11230.5
$ws.Range("M126").Value = -6661.25
$ws.Range("N126").Value = -16170.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1790.4706
$ws.Range("I16").Value = 1825.9
$ws.Range("K16").Value = 1825.9
$ws.Range("M16").Value = -1655.9
$ws.Range("H46").Value = 2209.75
$ws.Range("I46").Value = 1000
$ws.Range("J46").Value = 2613
$ws.Range("K46").Value = 1000
$ws.Range("L46").Value = 2613
$ws.Range("M46").Value = -812
$ws.Range("N46").Value = -2989
$ws.Range("H136").Value = 3144.5483
$ws.Range("I136").Value = 2687.682
$ws.Range("K136").Value = 8063.045999999999
$ws.Range("M136").Value = -5513.045999999999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 32142.857
$ws.Range("J123").Value = 32142.857
$ws.Range("L123").Value = 32142.857
$ws.Range("N123").Value = -41942.857
$ws.Range("H125").Value = 40000
$ws.Range("J125").Value = 40000
$ws.Range("L125").Value = 40000
$ws.Range("N125").Value = -49840
$ws.Range("H138").Value = 140000
$ws.Range("J138").Value = 140000
$ws.Range("L138").Value = 140000
$ws.Range("N138").Value = -150280
